# Rebuild the results dataframe dump on Sheet1 to reflect the new
# "First Estimate" / "Second Estimate" layout, broken out by age
# quartile bin and by individual preference question.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$statB = "45.82 `n(24.52) `np0.0 = 0.00 `np0.25 = 25.00 `np0.5 = 45.00 `np0.75 = 65.00 `np1.0 = 100.00 `nN=2980"
$statC = "0.41 `n(0.49) `n0: 0.59 `n1: 0.41 `nN=2980"

# ---------------------------------------------------------------
# Header block (rows 1-4) -- just re-label a couple of cells.
# ---------------------------------------------------------------
$ws.Range("A5").Value2 = "First Estimate"
$ws.Range("B4").Value2 = "Second Estimate"
$ws.Range("C4").Value2 = "% of Time Second Estimate is Preferred to First Estimate"

# ---------------------------------------------------------------
# Data for the four age bins (rows 6-9), re-using the formatting
# that already exists on row 6.
# ---------------------------------------------------------------
$bins = @("(-0.001, 25.0]", "(25.0, 47.0]", "(47.0, 65.0]", "(65.0, 100.0]")

for ($i = 0; $i -lt $bins.Length; $i++) {
    $r = 6 + $i
    if ($r -ne 6) {
        $ws.Range("A6:C6").Copy() | Out-Null
        $ws.Range("A" + $r + ":C" + $r).PasteSpecial($xlPasteFormats) | Out-Null
    }
    $ws.Range("A" + $r).Value2 = $bins[$i]
    $ws.Range("B" + $r).Value2 = $statB
    $ws.Range("C" + $r).Value2 = $statC
    $ws.Rows.Item($r).AutoFit() | Out-Null
}

# ---------------------------------------------------------------
# "Preference question" section header (row 11), formatted like
# the existing "Pooled" label in row 5.
# ---------------------------------------------------------------
$ws.Range("A5:A5").Copy() | Out-Null
$ws.Range("A11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A11").Value2 = "Preference question"

# ---------------------------------------------------------------
# Data for each individual preference question (rows 12-21).
# ---------------------------------------------------------------
$questions = @(
    " are more curious to try Soap flavored Bertie Bott's Every Flavour Beans than Earthworm flavored Bertie Bott's Every Flavour Beans",
    " believe climate change is a serious threat to human survival",
    " believe it should be legal to use psychedelic drugs  (e.g. LSD)",
    " can juggle with 3 balls",
    " own crypto currency (e.g. Bitcoin)",
    " prefer red over blue",
    " would feel safer in a self-driving car (rather than driving themselves)",
    " would prefer Dick Cheney over Sarah Palin to be president of the United States",
    " would rather have the superpower of flight than super strength",
    " would rather spend their evening going to a play at the theater than seeing a movie"
)

for ($i = 0; $i -lt $questions.Length; $i++) {
    $r = 12 + $i
    $ws.Range("A6:C6").Copy() | Out-Null
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("A" + $r).Value2 = $questions[$i]
    $ws.Range("B" + $r).Value2 = $statB
    $ws.Range("C" + $r).Value2 = $statC
    $ws.Rows.Item($r).AutoFit() | Out-Null
}

# ---------------------------------------------------------------
# Final "Pooled" section (rows 23-24).
# ---------------------------------------------------------------
$ws.Range("A5:A5").Copy() | Out-Null
$ws.Range("A23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A23").Value2 = "Pooled"

$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A24:C24").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A24").Value2 = "---"
$ws.Range("B24").Value2 = $statB
$ws.Range("C24").Value2 = $statC
$ws.Rows.Item(24).AutoFit() | Out-Null

$wb.Application.CutCopyMode = 0
